$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 176 (current "Musa paradisiaca" row), shifting
# everything from row 176 down through row 211 to rows 177-212.
$ws.Rows.Item(176).Insert()

# Populate the newly inserted row 176 with the "Monstera involuta" entry.
# Only columns A and E are populated (B, C, D left blank), matching the diff.
$ws.Cells.Item(176, 1).Value = "Monstera involuta"
$ws.Cells.Item(176, 5).Value = "shrub"
